$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Outcomes")

# Insert 3 new rows above the existing row 14 ("Persistent absentees ... 10%"),
# pushing all subsequent rows down by 3 (old row 14 -> new row 17, etc.)
$ws.Range("A14:A16").EntireRow.Insert()

# Row 14: Severe absentees for CINO at 31 March (overall absence 50% or more)
$ws.Range("A14").Value2 = 13
$ws.Range("B14").Value2 = "Outcomes"
$ws.Range("C14").Value2 = "Outcome 1: Children, young people and families stay together"
$ws.Range("D14").Value2 = "Child wellbeing and development"
$ws.Range("E14").Value2 = "Severe absentees for CINO at 31 March (overall absence 50% or more)"
$ws.Range("F14").Value2 = "outcomes_absence"
$ws.Range("G14").Value2 = "pt_pupils_pa_50_exact"
$ws.Range("H14").Value2 = "percent"
$ws.Range("I14").Value2 = "list('social_care_group' = 'CINO at 31 March', 'school_type' = 'Total')"

# Row 15: Severe absentees for CPPO at 31 March (overall absence 50% or more)
$ws.Range("A15").Value2 = 14
$ws.Range("B15").Value2 = "Outcomes"
$ws.Range("C15").Value2 = "Outcome 1: Children, young people and families stay together"
$ws.Range("D15").Value2 = "Child wellbeing and development"
$ws.Range("E15").Value2 = "Severe absentees for CPPO at 31 March (overall absence 50% or more)"
$ws.Range("F15").Value2 = "outcomes_absence"
$ws.Range("G15").Value2 = "pt_pupils_pa_50_exact"
$ws.Range("H15").Value2 = "percent"
$ws.Range("I15").Value2 = "list('social_care_group' = 'CPPO at 31 March', 'school_type' = 'Total')"

# Row 16: Severe absentees for CLA 12 months at 31 March (overall absence 50% or more)
$ws.Range("A16").Value2 = 15
$ws.Range("B16").Value2 = "Outcomes"
$ws.Range("C16").Value2 = "Outcome 1: Children, young people and families stay together"
$ws.Range("D16").Value2 = "Child wellbeing and development"
$ws.Range("E16").Value2 = "Severe absentees for CLA 12 months at 31 March (overall absence 50% or more)"
$ws.Range("F16").Value2 = "outcomes_absence"
$ws.Range("G16").Value2 = "pt_pupils_pa_50_exact"
$ws.Range("H16").Value2 = "percent"
$ws.Range("I16").Value2 = "list('social_care_group' = 'CLA 12 months at 31 March', 'school_type' = 'Total')"

# The values in column A are a plain sequential row index (not a formula),
# so renumber every pre-existing row below the insertion point to keep the
# A-column sequence (row number - 1) continuous after the insert.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 17; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 1
}

# Match the saved selection state recorded in the commit (cell H14 selected).
$ws.Range("H14").Select()
